$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Baboon)
$ws.Range("B2").Value = 93
$ws.Range("C2").Value = 0.00477897934615612
$ws.Range("D2").Value = 20
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 19
$ws.Range("G2").Value = -22
$ws.Range("H2").Value = 3
$ws.Range("I2").Value = 0.999989
$ws.Range("J2").Value = 0.000007
$ws.Range("K2").Value = 0.000004
$ws.Range("L2").Value = 0.999657
$ws.Range("M2").Value = 0.000343
$ws.Range("N2").Value = 0.999993
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0.000007
$ws.Range("R2").Value = 0.999996
$ws.Range("S2").Value = 0.999654
$ws.Range("T2").Value = 0.000007
$ws.Range("U2").Value = 0

# Row 3 (Hummingbird)
$ws.Range("B3").Value = 29
$ws.Range("C3").Value = 0.02764503099024296
$ws.Range("D3").Value = 4
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 4
$ws.Range("G3").Value = -5
$ws.Range("H3").Value = 1
$ws.Range("I3").Value = 0.965796
$ws.Range("J3").Value = -0.007907000000000001
$ws.Range("K3").Value = 0.016894
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = 0
$ws.Range("Q3").Value = 1
$ws.Range("R3").Value = 0.001104
$ws.Range("S3").Value = 0
$ws.Range("T3").Value = 0.999448
$ws.Range("U3").Value = 0.000552

# Row 4 (Panther)
$ws.Range("B4").Value = 110
$ws.Range("C4").Value = 0.005778418853878975
$ws.Range("D4").Value = 13
$ws.Range("E4").Value = 4
$ws.Range("F4").Value = 11
$ws.Range("G4").Value = -15
$ws.Range("H4").Value = 4
$ws.Range("I4").Value = 0.000004
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0.999996
$ws.Range("L4").Value = 0.000339
$ws.Range("M4").Value = 0.999661
$ws.Range("N4").Value = 0.999996
$ws.Range("O4").Value = 0.000007
$ws.Range("P4").Value = 0.000343
$ws.Range("Q4").Value = 0.000004
$ws.Range("R4").Value = 0.999989
$ws.Range("S4").Value = 0.999996
$ws.Range("T4").Value = 0
$ws.Range("U4").Value = 0.000004

